$d = $word.ActiveDocument

# The <id> tag text for this page was split across three runs:
#   "<id>p14" + "8" + "v_1</id> "
# Collapse them into a single run with the full text "<id>p148v_1</id> ".
$d.Content.Find.Execute("<id>p148v_1</id> ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p148v_1</id> ", 2)
